$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value2 = 111943940
$ws.Range("B3").Value2 = 77515
$ws.Range("D3").Value2 = "NT"
$ws.Range("E3").Value2 = 6425
$ws.Range("F3").Value2 = "Garnlav"
$ws.Range("G3").Value2 = "Alectoria sarmentosa"
$ws.Range("H3").Value2 = "(Ach.) Ach."
$ws.Range("Q3").Value2 = 600236.5842754361
$ws.Range("R3").Value2 = 7221446.606380152
$ws.Range("AX3").Value2 = "Maja Östlund, Simon Mattsson"

# Row 4
$ws.Range("A4").Value2 = 111943816
$ws.Range("B4").Value2 = 90682
$ws.Range("D4").Value2 = "NT"
$ws.Range("E4").Value2 = 2059
$ws.Range("F4").Value2 = "Skrovlig taggsvamp"
$ws.Range("G4").Value2 = "Hydnellum scabrosum"
$ws.Range("H4").Value2 = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q4").Value2 = 600428.2460409373
$ws.Range("R4").Value2 = 7221623.071005571
$ws.Range("AX4").Value2 = "Simon Mattsson, Maja Östlund"

# Row 6
$ws.Range("A6").Value2 = 111943907
$ws.Range("B6").Value2 = 90682
$ws.Range("D6").Value2 = "NT"
$ws.Range("E6").Value2 = 2059
$ws.Range("F6").Value2 = "Skrovlig taggsvamp"
$ws.Range("G6").Value2 = "Hydnellum scabrosum"
$ws.Range("H6").Value2 = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q6").Value2 = 600408.8168718456
$ws.Range("R6").Value2 = 7221648.141257811
$ws.Range("AX6").Value2 = "Maja Östlund, Simon Mattsson"

# Row 7
$ws.Range("A7").Value2 = 111943887
$ws.Range("B7").Value2 = 90660
$ws.Range("D7").Value2 = "NT"
$ws.Range("E7").Value2 = 4362
$ws.Range("F7").Value2 = "Blå taggsvamp"
$ws.Range("G7").Value2 = "Hydnellum caeruleum"
$ws.Range("H7").Value2 = "(Hornem.) P.Karst."
$ws.Range("Q7").Value2 = 600485.380207623
$ws.Range("R7").Value2 = 7221469.788948845
$ws.Range("AX7").Value2 = "Simon Mattsson, Maja Östlund"

# Row 8
$ws.Range("A8").Value2 = 111943881
$ws.Range("B8").Value2 = 90666
$ws.Range("D8").Value2 = "LC"
$ws.Range("E8").Value2 = 4364
$ws.Range("F8").Value2 = "Dropptaggsvamp"
$ws.Range("G8").Value2 = "Hydnellum ferrugineum"
$ws.Range("H8").Value2 = "(Fr.:Fr.) P. Karst."
$ws.Range("Q8").Value2 = 600419.1458058911
$ws.Range("R8").Value2 = 7221629.933000125
$ws.Range("AX8").Value2 = "Simon Mattsson, Maja Östlund"

# Row 9
$ws.Range("A9").Value2 = 111943877
$ws.Range("B9").Value2 = 90666
$ws.Range("D9").Value2 = "LC"
$ws.Range("E9").Value2 = 4364
$ws.Range("F9").Value2 = "Dropptaggsvamp"
$ws.Range("G9").Value2 = "Hydnellum ferrugineum"
$ws.Range("H9").Value2 = "(Fr.:Fr.) P. Karst."
$ws.Range("Q9").Value2 = 600475.944652258
$ws.Range("R9").Value2 = 7221498.994947547
$ws.Range("AX9").Value2 = "Simon Mattsson, Maja Östlund"

# Row 10
$ws.Range("A10").Value2 = 111943879
$ws.Range("B10").Value2 = 90666
$ws.Range("D10").Value2 = "LC"
$ws.Range("E10").Value2 = 4364
$ws.Range("F10").Value2 = "Dropptaggsvamp"
$ws.Range("G10").Value2 = "Hydnellum ferrugineum"
$ws.Range("H10").Value2 = "(Fr.:Fr.) P. Karst."
$ws.Range("Q10").Value2 = 600452.4405386611
$ws.Range("R10").Value2 = 7221544.595615291
$ws.Range("AX10").Value2 = "Simon Mattsson, Maja Östlund"

# Row 11
$ws.Range("A11").Value2 = 111943815
$ws.Range("B11").Value2 = 90682
$ws.Range("D11").Value2 = "NT"
$ws.Range("E11").Value2 = 2059
$ws.Range("F11").Value2 = "Skrovlig taggsvamp"
$ws.Range("G11").Value2 = "Hydnellum scabrosum"
$ws.Range("H11").Value2 = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q11").Value2 = 600430.1552666676
$ws.Range("R11").Value2 = 7221629.040516024
$ws.Range("AX11").Value2 = "Simon Mattsson, Maja Östlund"

# Row 13
$ws.Range("A13").Value2 = 111943880
$ws.Range("B13").Value2 = 90666
$ws.Range("D13").Value2 = "LC"
$ws.Range("E13").Value2 = 4364
$ws.Range("F13").Value2 = "Dropptaggsvamp"
$ws.Range("G13").Value2 = "Hydnellum ferrugineum"
$ws.Range("H13").Value2 = "(Fr.:Fr.) P. Karst."
$ws.Range("Q13").Value2 = 600437.3706075938
$ws.Range("R13").Value2 = 7221453.399649266
$ws.Range("AX13").Value2 = "Simon Mattsson, Maja Östlund"

# Row 14
$ws.Range("A14").Value2 = 111943814
$ws.Range("B14").Value2 = 90682
$ws.Range("D14").Value2 = "NT"
$ws.Range("E14").Value2 = 2059
$ws.Range("F14").Value2 = "Skrovlig taggsvamp"
$ws.Range("G14").Value2 = "Hydnellum scabrosum"
$ws.Range("H14").Value2 = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q14").Value2 = 600437.2972375784
$ws.Range("R14").Value2 = 7221630.12595264
$ws.Range("AX14").Value2 = "Simon Mattsson, Maja Östlund"

# Row 15
$ws.Range("A15").Value2 = 111943882
$ws.Range("B15").Value2 = 90666
$ws.Range("D15").Value2 = "LC"
$ws.Range("E15").Value2 = 4364
$ws.Range("F15").Value2 = "Dropptaggsvamp"
$ws.Range("G15").Value2 = "Hydnellum ferrugineum"
$ws.Range("H15").Value2 = "(Fr.:Fr.) P. Karst."
$ws.Range("Q15").Value2 = 600418.6659407767
$ws.Range("R15").Value2 = 7221432.097987156
$ws.Range("AX15").Value2 = "Simon Mattsson, Maja Östlund"

# Row 16
$ws.Range("A16").Value2 = 111943803
$ws.Range("B16").Value2 = 89405
$ws.Range("D16").Value2 = "NT"
$ws.Range("E16").Value2 = 1202
$ws.Range("F16").Value2 = "Ullticka"
$ws.Range("G16").Value2 = "Phellinidium ferrugineofuscum"
$ws.Range("H16").Value2 = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q16").Value2 = 600424.4900947324
$ws.Range("R16").Value2 = 7221684.100621465
$ws.Range("AX16").Value2 = "Simon Mattsson, Maja Östlund"

# Row 17
$ws.Range("A17").Value2 = 111943947
$ws.Range("B17").Value2 = 85715
$ws.Range("D17").Value2 = "NT"
$ws.Range("E17").Value2 = 510
$ws.Range("F17").Value2 = "Doftskinn"
$ws.Range("G17").Value2 = "Cystostereum murrayi"
$ws.Range("H17").Value2 = "(Berk. & M.A. Curtis.) Pouzar"
$ws.Range("Q17").Value2 = 600352.1009737813
$ws.Range("R17").Value2 = 7221402.427547264
$ws.Range("AX17").Value2 = "Maja Östlund, Simon Mattsson"

# Row 18
$ws.Range("A18").Value2 = 111943841
$ws.Range("B18").Value2 = 77268
$ws.Range("D18").Value2 = "NT"
$ws.Range("E18").Value2 = 228912
$ws.Range("F18").Value2 = "Mörk kolflarnlav"
$ws.Range("G18").Value2 = "Carbonicola myrmecina"
$ws.Range("H18").Value2 = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q18").Value2 = 600367.4767540093
$ws.Range("R18").Value2 = 7221297.494507908
$ws.Range("AX18").Value2 = "Simon Mattsson, Maja Östlund"
